$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Binomial")

# Update the binomial distribution parameters (n, p, x) for the new exercise
$ws.Range("C6").Value = 15
$ws.Range("C7").Value = 0.85
$ws.Range("C8").Value = 10

# Change the "Resultado" column formatting from the custom 0.0000% format
# to the built-in 0.00% percentage format
$ws.Range("H7:H11").NumberFormat = "0.00%"

# Leave the selection where the author left it when saving
$ws.Range("G16").Select() | Out-Null
